$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("aburðartegundir")

# Forsendubreyting - uppfært niturinnihald (N, dálkur H) fyrir valda áburðarflokka
$ws.Range("H4").Value = 0.45
$ws.Range("H8").Value = 2.76
$ws.Range("H9").Value = 3.58
$ws.Range("H10").Value = 9.35
$ws.Range("H11").Value = 0.27
$ws.Range("H12").Value = 1.31
$ws.Range("H14").Value = 0.88

$ws.Activate()
$ws.Range("H15").Select()
